$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.842.09"
$ws.Range("E2").Value = "  +0.40%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.636.09"
$ws.Range("E3").Value = "  +0.62%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.16"
$ws.Range("E5").Value = "  +0.16%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5083"
$ws.Range("E6").Value = "  -0.29%  "

# Row 7
$ws.Range("E7").Value = "  +0.14%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2583"
$ws.Range("E8").Value = "  +0.94%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06426"
$ws.Range("E9").Value = "  +1.79%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.38"
$ws.Range("E10").Value = "  +5.29%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07799"
$ws.Range("E11").Value = "  +0.31%  "

# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.262"
$ws.Range("E12").Value = "  +0.93%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.646.36"
$ws.Range("E13").Value = "  +1.16%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.862.67"
$ws.Range("E14").Value = "  +0.64%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5600"
$ws.Range("E15").Value = "  +1.72%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7657"
$ws.Range("E16").Value = "  +2.47%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.24"
$ws.Range("E17").Value = "  -0.39%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.848.18"

# Row 19
$ws.Range("E19").Value = "  +0.12%  "

# Row 20
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.387"
$ws.Range("E20").Value = "  -0.25%  "

# Row 21
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.21"
$ws.Range("E21").Value = "  -0.15%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.985"
$ws.Range("E22").Value = "  +2.00%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.145"
$ws.Range("E23").Value = "  +2.71%  "

# Row 24
$ws.Range("E24").Value = "  +0.13%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.761"
$ws.Range("E25").Value = "  -6.57%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "138.53"
$ws.Range("E26").Value = "  -2.37%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1235"
$ws.Range("E27").Value = "  -1.52%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.849"
$ws.Range("E28").Value = "  +1.92%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.53"
$ws.Range("E29").Value = "  +0.35%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.241"
$ws.Range("E30").Value = "  +0.40%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04965"
$ws.Range("E31").Value = "  +1.93%  "

# Row 32
$ws.Range("E32").Value = "  +2.40%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.254"
$ws.Range("E33").Value = "  +3.39%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.568"
$ws.Range("E34").Value = "  +2.21%  "

# Row 35
$ws.Range("E35").Value = "  +0.52%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9033"
$ws.Range("E36").Value = "  +1.45%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.579"
$ws.Range("E37").Value = "  +1.56%  "

# Row 38
$ws.Range("E38").Value = "  +1.41%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.134.05"
$ws.Range("E39").Value = "  +2.13%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01568"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9964"
$ws.Range("E41").Value = "  -0.41%  "

# Row 42
$ws.Range("E42").Value = "  +2.25%  "

# Row 43
$ws.Range("E43").Value = "  -1.53%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7999"
$ws.Range("E44").Value = "  +0.43%  "

# Row 45
$ws.Range("E45").Value = "  -0.62%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.61"
$ws.Range("E46").Value = "  +2.03%  "

# Row 47
$ws.Range("E47").Value = "  -3.74%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.814"
$ws.Range("E48").Value = "  +4.06%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05031"
$ws.Range("E49").Value = "  -2.01%  "

# Row 50
$ws.Range("E50").Value = "  +0.07%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.003"
$ws.Range("E51").Value = "  +0.30%  "
